# Updates cryptos list D (Price) / E (Volume 1h) columns per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = new price text (optional); E = new volume text }
$updates = @(
    @{ Row = 2; D = '70.927.78'; E = '  +2.19%  ' }
    @{ Row = 3; D = '3.552.36'; E = '  +0.90%  ' }
    @{ Row = 4; D = $null; E = '  +0.04%  ' }
    @{ Row = 5; D = '607.32'; E = '  +4.13%  ' }
    @{ Row = 6; D = '175.03'; E = '  +0.58%  ' }
    @{ Row = 7; D = '0.618'; E = '  -0.61%  ' }
    @{ Row = 8; D = '3.547.44'; E = '  +1.02%  ' }
    @{ Row = 9; D = $null; E = '  -0.03%  ' }
    @{ Row = 10; D = '0.202'; E = '  +6.11%  ' }
    @{ Row = 11; D = '6.74'; E = '  -0.63%  ' }
    @{ Row = 12; D = '0.589'; E = '  -1.38%  ' }
    @{ Row = 13; D = '47.76'; E = '  +1.41%  ' }
    @{ Row = 14; D = $null; E = '  +1.54%  ' }
    @{ Row = 15; D = '4.125.22'; E = '  +1.15%  ' }
    @{ Row = 16; D = '629.86'; E = '  -7.08%  ' }
    @{ Row = 17; D = '8.47'; E = '  -3.49%  ' }
    @{ Row = 18; D = '70.961.89'; E = '  +2.32%  ' }
    @{ Row = 19; D = '3.552.00'; E = '  +0.91%  ' }
    @{ Row = 20; D = $null; E = '  -1.85%  ' }
    @{ Row = 21; D = '17.48'; E = '  -0.08%  ' }
    @{ Row = 22; D = $null; E = '  -10.03%  ' }
    @{ Row = 23; D = '0.892'; E = '  -1.52%  ' }
    @{ Row = 24; D = '15.97'; E = '  -1.27%  ' }
    @{ Row = 25; D = '97.25'; E = '  -0.94%  ' }
    @{ Row = 26; D = $null; E = '  -0.22%  ' }
    @{ Row = 27; D = $null; E = '  -0.09%  ' }
    @{ Row = 28; D = $null; E = '  -1.48%  ' }
    @{ Row = 29; D = '9.27'; E = '  -1.91%  ' }
    @{ Row = 30; D = '33.50'; E = '  +1.21%  ' }
    @{ Row = 31; D = '3.15'; E = '  -1.90%  ' }
    @{ Row = 32; D = $null; E = '  -3.13%  ' }
    @{ Row = 33; D = $null; E = '  -0.89%  ' }
    @{ Row = 34; D = $null; E = '  -3.38%  ' }
    @{ Row = 35; D = '569.75'; E = '  -4.49%  ' }
    @{ Row = 36; D = $null; E = '  +1.46%  ' }
    @{ Row = 37; D = '10.83'; E = '  -0.91%  ' }
    @{ Row = 38; D = $null; E = '  -2.04%  ' }
    @{ Row = 39; D = '57.58'; E = '  +0.54%  ' }
    @{ Row = 40; D = $null; E = '  -0.02%  ' }
    @{ Row = 41; D = '0.144'; E = '  +5.88%  ' }
    @{ Row = 42; D = $null; E = '  +3.16%  ' }
    @{ Row = 43; D = '0.330'; E = '  -2.03%  ' }
    @{ Row = 44; D = '3.351.00'; E = '  -2.03%  ' }
    @{ Row = 45; D = '3.05'; E = '  +4.51%  ' }
    @{ Row = 46; D = '0.0₃0724'; E = '  +1.82%  ' }
    @{ Row = 47; D = '33.23'; E = '  -0.64%  ' }
    @{ Row = 48; D = '2.67'; E = '  +2.26%  ' }
    @{ Row = 49; D = $null; E = '  -2.22%  ' }
    @{ Row = 50; D = '134.52'; E = '  +1.73%  ' }
    @{ Row = 51; D = '5.71'; E = '  -2.25%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Source values are plain text (e.g. "607.32"); Excel would otherwise
        # auto-convert a bare decimal into a Number on assignment, so force the
        # cell to Text, write the string, then drop the format change again so
        # the cell ends up back on the default (unstyled) style.
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
